$wb = $excel.ActiveWorkbook

# --- Work on "Repayment schedule" sheet: insert a new (blank) column before N ---
# This shifts the old N/O/P columns (Late / heading / Outstanding) one place
# right to O/P/Q, matching the "Variable Instalments" layout used elsewhere.
$ws = $wb.Worksheets.Item("Repayment schedule")
$ws.Columns("N").Insert()

# The inherited column width should match column M (its left neighbour),
# which is what Excel does by default when a column is inserted.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Activate this sheet (becomes the selected tab) and park the selection at S7
$ws.Activate()
$ws.Range("S7").Select()
